# "Footer Added & Survey Score Removed From DB/XLSX"
#
# The "Average Survey Score" column (column E) is no longer collected, so
# it is removed from the sheet entirely: the column is deleted (shifting
# Counter Queries Taken / Chats Taken / Breached Tickets / Total one slot
# to the left), and the "Total" formula no longer adds the old survey-score
# bonus term.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Average Survey Score" column -----------------------------
$ws.Range("E1").EntireColumn.Delete() | Out-Null

# --- A few underlying counts were corrected at the same time --------------
$ws.Range("F4").Value = 7      # Jenny  - Chats Taken
$ws.Range("G6").Value = 3      # Josh   - Breached Tickets
$ws.Range("F7").Value = 4      # Karina - Chats Taken

# --- Total no longer includes the old "(SurveyScore * 10)" bonus ----------
$ws.Range("H2").Formula = "=SUM(C2:G2)-(G2*50)"
$ws.Range("H3:H10").FormulaR1C1 = "=SUM(RC[-5]:RC[-1])-(RC[-1]*50)"

# Refresh the RANK() helper column now that it is keyed off the new Total
# column (H instead of I).
$ws.Range("A2:A10").FormulaR1C1 = "=RANK(RC[7],R2C8:R10C8)"

# --- Defined names shift along with the deleted column --------------------
$wb.Names("Data").RefersTo = "='Monthly Results'!`$B`$2:`$H`$10"
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Monthly Results'!`$J`$1:`$L`$2"
    }
}

# --- Selection left where the editor's cursor ended up ---------------------
$ws.Range("F6").Select() | Out-Null

$excel.Calculate()
